$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 - "Điểm danh": mark as "Chưa cần" priority, add note
$ws.Range("C40").Value = "Chưa cần"
$ws.Range("F40").Value = "Tạm chưa làm vì chưa cần lắm!!"

# Row 43 - "Sửa giờ chấm công": progress updated to backend done
$ws.Range("E43").Value = "xong front-end-Xong Backend"

# Row 44 - "Xem lịch sử sửa giờ": mark as "Chưa cần" priority, add note
$ws.Range("C44").Value = "Chưa cần"
$ws.Range("F44").Value = "Tạm chưa làm vì chưa cần lắm!!"

# Row 46 - "ATIN Smartface About": clear the person in charge (D46 was "Sáng")
$ws.Range("D46").Value = ""

# Row 47 - "Main Form": reassign person in charge, update progress note and add remark
$ws.Range("D47").Value = "xxxxxxxxxxxxxxxx"
$ws.Range("E47").Value = "Chưa đụng đến"
$ws.Range("F47").Value = "sau khi hoàn thiện các chức năng thì tiến hành merge"

# Move/select the active cell to F47, matching the final cursor position
$ws.Range("F47").Select()
